# daily auto push: 2026-01-05 18:48 UTC
# Insert a new row of data at row 564 ("2026/01/05", 月, 23, 19),
# pushing the existing rows 564-605 down to 565-606.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("564:564").Insert()

# Keep the date column as literal text (matches the rest of the sheet,
# which stores dates as strings rather than real date values) - without
# this, Excel auto-converts the "2026/01/05"-shaped string into a date
# serial number.
$ws.Range("A564").NumberFormat = "@"
$ws.Range("A564").Value = "2026/01/05"
$ws.Range("B564").Value = "月"
$ws.Range("C564").Value = 23
$ws.Range("D564").Value = 19
